$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This data dictionary worksheet underwent a re-sort of several variable
# blocks (rows identified by the "Variable" key in column T). Rather than
# reconstructing row order with physical cut/insert operations, we update
# every touched cell in place so each row ends up holding the field values
# of its new logical position -- equivalent to the net effect of the sort.

# Rows 20-22: bring the Age45_49 / Age50_54 / AgeOv18 trio (previously rows 32-34) to the top of the Social block
$ws.Range("C20").Value = 'x'
$ws.Range("D20").Value = 'x'
$ws.Range("E20").Value = 'x'
$ws.Range("K20").Value = 'x'
$ws.Range("M20").Value = ""
$ws.Range("P20").Value = ""
$ws.Range("T20").Value = 'Age45_49'
$ws.Range("U20").Value = 'Population: Age 45-49'
$ws.Range("V20").Value = 'Total population between age 45-49'
$ws.Range("X20").Value = 'ACS 2018, 5-Year; 2010 Decennial Census; IPUMS NHGIS'
$ws.Range("Y20").Value = 'American Community Survey 2014-2018 5 Year Estimates; 2010 Decennial Census; Integrated Public Use Microdata Service National Historic Geographic Information Systems'
$ws.Range("AA20").Value = 'integer'
$ws.Range("AB20").Value = '''467768'
$ws.Range("AD20").Value = '1990 and 2000 data from respective decennial censuses downloaded from IPUMS NHGIS and aggregated upwards.'
$ws.Range("C21").Value = 'x'
$ws.Range("D21").Value = 'x'
$ws.Range("E21").Value = 'x'
$ws.Range("K21").Value = 'x'
$ws.Range("M21").Value = ""
$ws.Range("P21").Value = ""
$ws.Range("T21").Value = 'Age50_54'
$ws.Range("U21").Value = 'Population: Age 50-54'
$ws.Range("V21").Value = 'Total population between age 50-54'
$ws.Range("X21").Value = 'ACS 2018, 5-Year; 2010 Decennial Census; IPUMS NHGIS'
$ws.Range("Y21").Value = 'American Community Survey 2014-2018 5 Year Estimates; 2010 Decennial Census; Integrated Public Use Microdata Service National Historic Geographic Information Systems'
$ws.Range("AB21").Value = '''476486'
$ws.Range("AD21").Value = '1990 and 2000 data from respective decennial censuses downloaded from IPUMS NHGIS and aggregated upwards.'
$ws.Range("E22").Value = 'x'
$ws.Range("K22").Value = 'x'
$ws.Range("M22").Value = ""
$ws.Range("P22").Value = ""
$ws.Range("T22").Value = 'AgeOv18'
$ws.Range("U22").Value = 'Population: Age 18+'
$ws.Range("V22").Value = 'Total population at or over age 18'
$ws.Range("X22").Value = 'ACS 2018, 5-Year; 2010 Decennial Census'
$ws.Range("Y22").Value = 'American Community Survey 2014-2018 5 Year Estimates; 2010 Decennial Census'
$ws.Range("AB22").Value = '''5661461'
# Rows 23-34: the remaining Social/demographic rows shift down by three to make room
$ws.Range("T23").Value = 'FemP'
$ws.Range("U23").Value = '% Population that is Female'
$ws.Range("T24").Value = 'MaleP'
$ws.Range("U24").Value = '% Population that is Male'
$ws.Range("T25").Value = 'MedAge'
$ws.Range("U25").Value = 'Median age'
$ws.Range("T26").Value = 'Ovr16P'
$ws.Range("U26").Value = '% Population over 16 years'
$ws.Range("T27").Value = 'Ovr18P'
$ws.Range("U27").Value = '% Population over 18 years'
$ws.Range("T28").Value = 'Ovr21P'
$ws.Range("U28").Value = '% Population over 21 years'
$ws.Range("T29").Value = 'Ovr62P'
$ws.Range("U29").Value = '% Population over 62 years'
$ws.Range("T30").Value = 'SRatio'
$ws.Range("U30").Value = 'Sex ratio for the total population (males per 100 females)'
$ws.Range("T31").Value = 'SRatio18'
$ws.Range("U31").Value = 'Sex ratio among adults aged 18 and older (males per 100 females)'
$ws.Range("C32").Value = ""
$ws.Range("D32").Value = ""
$ws.Range("E32").Value = ""
$ws.Range("K32").Value = ""
$ws.Range("M32").Value = 'x'
$ws.Range("P32").Value = 'x'
$ws.Range("T32").Value = 'SRatio65'
$ws.Range("U32").Value = 'Sex ratio among seniors aged 65 and older (males per 100 females)'
$ws.Range("V32").Value = ""
$ws.Range("X32").Value = 'ACS '
$ws.Range("Y32").Value = 'American Community Survey (5-Year Estimate)'
$ws.Range("AA32").Value = 'number'
$ws.Range("AB32").Value = '''60.0'
$ws.Range("AD32").Value = ""
$ws.Range("C33").Value = ""
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = ""
$ws.Range("K33").Value = ""
$ws.Range("M33").Value = 'x'
$ws.Range("P33").Value = 'x'
$ws.Range("T33").Value = 'Und18P'
$ws.Range("U33").Value = '% Population under 18 years old'
$ws.Range("V33").Value = ""
$ws.Range("X33").Value = 'ACS '
$ws.Range("Y33").Value = 'American Community Survey (5-Year Estimate)'
$ws.Range("AB33").Value = '''60.0'
$ws.Range("AD33").Value = ""
$ws.Range("E34").Value = ""
$ws.Range("K34").Value = ""
$ws.Range("M34").Value = 'x'
$ws.Range("P34").Value = 'x'
$ws.Range("T34").Value = 'Und5P'
$ws.Range("U34").Value = '% Population under 5 years old '
$ws.Range("V34").Value = ""
$ws.Range("X34").Value = 'ACS '
$ws.Range("Y34").Value = 'American Community Survey (5-Year Estimate)'
$ws.Range("AB34").Value = '''60.0'
# Rows 79-117: 'Nal...' naloxone variable codes corrected to 'Nalt...'
$ws.Range("T79").Value = 'NaltRm30'
$ws.Range("T80").Value = 'NaltRm60'
$ws.Range("T81").Value = 'NaltRm90'
$ws.Range("T95").Value = 'NaltMinDis'
$ws.Range("T96").Value = 'NaltTmDr'
$ws.Range("T97").Value = 'NaltCntDr30'
$ws.Range("T104").Value = 'NaltTmWk'
$ws.Range("T105").Value = 'NaltCntWk60'
$ws.Range("T106").Value = 'NaltCntWk30'
$ws.Range("T113").Value = 'NaltTmBk'
$ws.Range("T114").Value = 'NaltCntBk60'
$ws.Range("T115").Value = 'NaltCntBk30'
$ws.Range("T117").Value = 'NaltCntDr60'
# Rows 122-124: OTP access rows rotate (OtpMinDis -> bottom, OtpTmDr/OtpCntDr move up)
$ws.Range("N122").Value = 'x'
$ws.Range("R122").Value = 'x'
$ws.Range("T122").Value = 'OtpTmDr'
$ws.Range("U122").Value = 'Driving Time (min) to nearest Opioid Treatment Program (OTP)'
$ws.Range("V122").Value = 'Driving time from tract/zip origin centroid to the nearest tract/zip OTP destination centroid, in minutes'
$ws.Range("X122").Value = 'SAMSHA, 2021'
$ws.Range("Y122").Value = 'U.S. Substance Abuse and Mental Health Services Administration Treatment Locator, 2021'
$ws.Range("AB122").Value = '''27.39'
$ws.Range("AD122").Value = 'This dataset includes all US states, Washington D.C., and Puerto Rico. It does not include the territories Guam, Northern Mariana Islands, American Samoa, Palau. Zip code and tract centroids are not population-weighted.'
$ws.Range("R123").Value = ""
$ws.Range("T123").Value = 'OtpCntDr'
$ws.Range("U123").Value = 'Count of Opioid Treatment Programs (OTP) (30-min drive)'
$ws.Range("V123").Value = 'Count of OTPs within a 30-minute driving threshold'
$ws.Range("AA123").Value = 'integer'
$ws.Range("AB123").Value = '''1'
$ws.Range("N124").Value = ""
$ws.Range("T124").Value = 'OtpMinDis'
$ws.Range("U124").Value = 'Distance to nearest OTP'
$ws.Range("V124").Value = 'Euclidean distance* from tract/zip centroid to nearest OTP service location, in miles'
$ws.Range("X124").Value = 'SAMHSA'
$ws.Range("Y124").Value = 'The Substance Abuse and Mental Health Services Administration'
$ws.Range("AA124").Value = 'number'
$ws.Range("AB124").Value = '''121.0'
$ws.Range("AD124").Value = ""
# Rows 164-165: swap Essential Workers count/percent rows
$ws.Range("M164").Value = ""
$ws.Range("P164").Value = ""
$ws.Range("T164").Value = 'EssnWrkE'
$ws.Range("U164").Value = 'Count of Essential Workers'
$ws.Range("V164").Value = 'Estimated count of population employed in essential occupations.'
$ws.Range("AA164").Value = 'integer'
$ws.Range("AB164").Value = '''1509709'
$ws.Range("M165").Value = 'x'
$ws.Range("P165").Value = 'x'
$ws.Range("T165").Value = 'EssnWrkP'
$ws.Range("U165").Value = 'Essential Workers %'
$ws.Range("V165").Value = 'Percentage of population employed in essential occupations.'
$ws.Range("AA165").Value = 'number'
$ws.Range("AB165").Value = '''42.96'
# Rows 166-170: SDOH index rows rotate (NeighbTyp moves to the bottom)
$ws.Range("G166").Value = 'x'
$ws.Range("M166").Value = ""
$ws.Range("R166").Value = 'x'
$ws.Range("T166").Value = 'SocEcAdvIn'
$ws.Range("U166").Value = 'Socioeconomic Advantage Index'
$ws.Range("V166").Value = 'Raw Socioeconomic Advantage Index (https://sdohatlas.github.io/)'
$ws.Range("X166").Value = 'GeoDa Data and Lab; SDOH Atlas'
$ws.Range("Y166").Value = 'GeoDa Data and Lab; Spatial Deterimants of Health Atlas'
$ws.Range("AA166").Value = 'number'
$ws.Range("AB166").Value = '''1.17'
$ws.Range("AC166").Value = ""
$ws.Range("T167").Value = 'LimMobInd'
$ws.Range("U167").Value = 'Limited Moblility Index'
$ws.Range("V167").Value = 'Raw Limited Mobility Index (https://sdohatlas.github.io/)'
$ws.Range("AB167").Value = '''0.54'
$ws.Range("T168").Value = 'UrbCoreInd'
$ws.Range("U168").Value = 'Urban Core Opportunity Index'
$ws.Range("V168").Value = 'Raw Urban Core Opportunity Index (https://sdohatlas.github.io/)'
$ws.Range("AB168").Value = '''-0.06'
$ws.Range("T169").Value = 'MicaInd'
$ws.Range("U169").Value = 'Mixed Immigrant Cohesion and Accesibility (MICA) Index'
$ws.Range("V169").Value = 'Raw Mixed Immigrant Cohesion and Accessibility (MICA) Index (https://sdohatlas.github.io/)'
$ws.Range("AB169").Value = '''0.59'
$ws.Range("G170").Value = ""
$ws.Range("M170").Value = 'x'
$ws.Range("R170").Value = ""
$ws.Range("T170").Value = 'NeighbTyp'
$ws.Range("U170").Value = 'Neighborhood Type'
$ws.Range("V170").Value = 'Categorical, one of seven neighborhood (tract-level) typologies: 1 = Rural Affordable; 2 = Suburban Affluent; 3 = Suburban Affordable; 4 = Extreme Poverty; 5 = Multilingual Working; 6 = Urban Core Opportunity; 7 = Sparse Areas'
$ws.Range("X170").Value = 'Kolak et al., 2020'
$ws.Range("Y170").Value = 'Kolak et al., 2020'
$ws.Range("AA170").Value = 'string'
$ws.Range("AB170").Value = '''3'
$ws.Range("AC170").Value = 'Some tracts are not assigned with any typologies because data are missing for factor analysis.'
